$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths between column A and column B
$ws.Columns.Item(1).ColumnWidth = 2.140625
$ws.Columns.Item(2).ColumnWidth = 3.140625

# Swap the values in columns A and B for rows 3-10
for ($r = 3; $r -le 10; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value()
    $bVal = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 1).Value = $bVal
    $ws.Cells.Item($r, 2).Value = $aVal
}
